$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 746
$ws.Range("K3").Value = 707
$ws.Range("E4").Value = 2020
$ws.Range("F4").Value = 1906
$ws.Range("J4").Value = 1792
$ws.Range("K4").Value = 155
$ws.Range("K5").Value = 43
$ws.Range("K6").Value = 997
$ws.Range("E7").Value = 26025
$ws.Range("F7").Value = 24099
$ws.Range("J7").Value = 29245
$ws.Range("K7").Value = 2648

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 42
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 33
$ws.Range("K3").Value = 40
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 17
$ws.Range("K3").Value = 27
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 21
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 76
$ws.Range("K8").Value = 166
$ws.Range("K9").Value = 13
$ws.Range("K17").Value = 4
$ws.Range("K18").Value = 17
$ws.Range("K19").Value = 67
$ws.Range("K22").Value = 9
$ws.Range("K27").Value = 33
$ws.Range("K29").Value = 130
$ws.Range("K33").Value = 121
$ws.Range("K37").Value = 80
$ws.Range("K41").Value = 25
$ws.Range("K42").Value = 89
$ws.Range("J43").Value = 249
$ws.Range("K49").Value = 20
$ws.Range("K52").Value = 69
$ws.Range("K53").Value = 30
$ws.Range("K54").Value = 50
$ws.Range("K55").Value = 25
$ws.Range("K60").Value = 20
$ws.Range("K61").Value = 3
$ws.Range("E63").Value = 360
$ws.Range("F63").Value = 192
$ws.Range("J63").Value = 86
$ws.Range("K63").Value = 11
$ws.Range("K64").Value = 13
$ws.Range("K65").Value = 70
$ws.Range("K67").Value = 111
$ws.Range("K73").Value = 27
$ws.Range("K79").Value = 68
$ws.Range("K83").Value = 49
$ws.Range("K85").Value = 129
$ws.Range("K86").Value = 20
$ws.Range("K88").Value = 36
$ws.Range("K91").Value = 27
$ws.Range("J94").Value = 331
$ws.Range("K94").Value = 32
$ws.Range("K96").Value = 42
$ws.Range("K97").Value = 21
$ws.Range("E101").Value = 26025
$ws.Range("F101").Value = 24099
$ws.Range("J101").Value = 29245
$ws.Range("K101").Value = 2648

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 2
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 38
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 18
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 19
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 35
$ws.Range("K3").Value = 20
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 11
$ws.Range("J4").Value = 27
$ws.Range("J7").Value = 331
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 6
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 1
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 4
$ws.Range("K3").Value = 1

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 23
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 18
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K2").Value = 1
$ws.Range("K7").Value = 3
